$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.452596788596793
$ws.Cells.Item(2, 3).Value = 0.2185145910852953
$ws.Cells.Item(2, 5).Value = 0.1708062654801981
$ws.Cells.Item(2, 6).Value = 2.491899588884806
$ws.Cells.Item(2, 7).Value = 1.457257625380507
$ws.Cells.Item(2, 8).Value = 1.320026110066692
$ws.Cells.Item(2, 10).Value = 0.08445800516152602
$ws.Cells.Item(2, 12).Value = 0.358910530734839
$ws.Cells.Item(2, 13).Value = 0.3596580683717008
$ws.Cells.Item(3, 2).Value = 1.361541582130997
$ws.Cells.Item(3, 3).Value = 0.2039858855056309
$ws.Cells.Item(3, 5).Value = 0.1712053721254572
$ws.Cells.Item(3, 6).Value = 2.490324633412328
$ws.Cells.Item(3, 7).Value = 1.451293054058141
$ws.Cells.Item(3, 8).Value = 1.323837418399052
$ws.Cells.Item(3, 10).Value = 0.08510190523331218
$ws.Cells.Item(3, 12).Value = 0.3556043110205707
$ws.Cells.Item(3, 13).Value = 0.3453601590858
$ws.Cells.Item(4, 2).Value = 1.306211074710745
$ws.Cells.Item(4, 3).Value = 0.1949808204091852
$ws.Cells.Item(4, 5).Value = 0.1714702012093472
$ws.Cells.Item(4, 6).Value = 2.490710674386747
$ws.Cells.Item(4, 7).Value = 1.448674861269083
$ws.Cells.Item(4, 8).Value = 1.326899498435949
$ws.Cells.Item(4, 10).Value = 0.08551877028405563
$ws.Cells.Item(4, 12).Value = 0.3536867675908866
$ws.Cells.Item(4, 13).Value = 0.3367297139114527
$ws.Cells.Item(5, 2).Value = 1.283809354036407
$ws.Cells.Item(5, 3).Value = 0.1912899413631095
$ws.Cells.Item(5, 5).Value = 0.171583105896739
$ws.Cells.Item(5, 6).Value = 2.491207879623985
$ws.Cells.Item(5, 7).Value = 1.447869558273595
$ws.Cells.Item(5, 8).Value = 1.328328573201588
$ws.Cells.Item(5, 10).Value = 0.08569406562051007
$ws.Cells.Item(5, 12).Value = 0.3529337698535144
$ws.Cells.Item(5, 13).Value = 0.3332502565166848
$ws.Cells.Item(6, 2).Value = 1.280098397017639
$ws.Cells.Item(6, 3).Value = 0.1906757901354155
$ws.Cells.Item(6, 5).Value = 0.1716021550474434
$ws.Cells.Item(6, 6).Value = 2.491310958532296
$ws.Cells.Item(6, 7).Value = 1.447751613626011
$ws.Cells.Item(6, 8).Value = 1.328576808980316
$ws.Cells.Item(6, 10).Value = 0.08572350093344205
$ws.Cells.Item(6, 12).Value = 0.3528104553242954
$ws.Cells.Item(6, 13).Value = 0.3326747664216185
$ws.Cells.Item(7, 2).Value = 1.30590836563988
$ws.Cells.Item(7, 3).Value = 0.1949311299576522
$ws.Cells.Item(7, 5).Value = 0.1714717036823057
$ws.Cells.Item(7, 6).Value = 2.490716004140552
$ws.Cells.Item(7, 7).Value = 1.448662942547145
$ws.Cells.Item(7, 8).Value = 1.32691803793567
$ws.Cells.Item(7, 10).Value = 0.08552111242059013
$ws.Cells.Item(7, 12).Value = 0.3536764971663331
$ws.Cells.Item(7, 13).Value = 0.3366826365905311
$ws.Cells.Item(8, 2).Value = 1.421081379734915
$ws.Cells.Item(8, 3).Value = 0.2135225877408971
$ws.Cells.Item(8, 5).Value = 0.1709397800299191
$ws.Cells.Item(8, 6).Value = 2.491075496575334
$ws.Cells.Item(8, 7).Value = 1.454983789236692
$ws.Cells.Item(8, 8).Value = 1.321190242994035
$ws.Cells.Item(8, 10).Value = 0.08467556456568737
$ws.Cells.Item(8, 12).Value = 0.3577472753629749
$ws.Cells.Item(8, 13).Value = 0.3546974179439033
$ws.Cells.Item(9, 2).Value = 1.65150570309055
$ws.Cells.Item(9, 3).Value = 0.2493136593951988
$ws.Cells.Item(9, 5).Value = 0.170053074464676
$ws.Cells.Item(9, 6).Value = 2.502535601017627
$ws.Cells.Item(9, 7).Value = 1.475706948247648
$ws.Cells.Item(9, 8).Value = 1.315699995413581
$ws.Cells.Item(9, 10).Value = 0.0831876094287658
$ws.Cells.Item(9, 12).Value = 0.3666176018837604
$ws.Cells.Item(9, 13).Value = 0.3911976810906168
$ws.Cells.Item(10, 2).Value = 1.823586430022829
$ws.Cells.Item(10, 3).Value = 0.2752093030404978
$ws.Cells.Item(10, 5).Value = 0.1694962607686037
$ws.Cells.Item(10, 6).Value = 2.517543306247248
$ws.Cells.Item(10, 7).Value = 1.496073565914173
$ws.Cells.Item(10, 8).Value = 1.315187564430687
$ws.Cells.Item(10, 10).Value = 0.0821974747890053
$ws.Cells.Item(10, 12).Value = 0.3736700483857618
$ws.Cells.Item(10, 13).Value = 0.4187261621301843
$ws.Cells.Item(11, 2).Value = 1.902477798061284
$ws.Cells.Item(11, 3).Value = 0.2869046004107076
$ws.Cells.Item(11, 5).Value = 0.1692633597001567
$ws.Cells.Item(11, 6).Value = 2.525808315687556
$ws.Cells.Item(11, 7).Value = 1.506469363458166
$ws.Cells.Item(11, 8).Value = 1.31572367769212
$ws.Cells.Item(11, 10).Value = 0.08176928369051062
$ws.Cells.Item(11, 12).Value = 0.3769935432811735
$ws.Cells.Item(11, 13).Value = 0.4314036501201812
$ws.Cells.Item(12, 2).Value = 1.9324395318651
$ws.Cells.Item(12, 3).Value = 0.2913211933253876
$ws.Cells.Item(12, 5).Value = 0.1691780877956974
$ws.Cells.Item(12, 6).Value = 2.529145323188018
$ws.Cells.Item(12, 7).Value = 1.510569676906584
$ws.Cells.Item(12, 8).Value = 1.316037654661756
$ws.Cells.Item(12, 10).Value = 0.08161032605152485
$ws.Cells.Item(12, 12).Value = 0.3782685338278355
$ws.Cells.Item(12, 13).Value = 0.4362264173262034
$ws.Cells.Item(13, 2).Value = 1.925982863281149
$ws.Cells.Item(13, 3).Value = 0.2903705414960029
$ws.Cells.Item(13, 5).Value = 0.1691963228014561
$ws.Cells.Item(13, 6).Value = 2.528417415125318
$ws.Cells.Item(13, 7).Value = 1.509679306284227
$ws.Cells.Item(13, 8).Value = 1.315965093426627
$ws.Cells.Item(13, 10).Value = 0.08164441869617889
$ws.Cells.Item(13, 12).Value = 0.3779932119285405
$ws.Cells.Item(13, 13).Value = 0.4351867687139261
$ws.Cells.Item(14, 2).Value = 1.904941021820207
$ws.Cells.Item(14, 3).Value = 0.2872681997590973
$ws.Cells.Item(14, 5).Value = 0.1692562858127467
$ws.Cells.Item(14, 6).Value = 2.526078697533052
$ws.Cells.Item(14, 7).Value = 1.506803413031378
$ws.Cells.Item(14, 8).Value = 1.315747282903715
$ws.Cells.Item(14, 10).Value = 0.08175614227539363
$ws.Cells.Item(14, 12).Value = 0.3770981084454661
$ws.Cells.Item(14, 13).Value = 0.4317999805141142
$ws.Cells.Item(15, 2).Value = 1.892063636781813
$ws.Cells.Item(15, 3).Value = 0.2853663427933384
$ws.Cells.Item(15, 5).Value = 0.1692933952048143
$ws.Cells.Item(15, 6).Value = 2.524673166208501
$ws.Cells.Item(15, 7).Value = 1.50506318797747
$ws.Cells.Item(15, 8).Value = 1.315628328670471
$ws.Cells.Item(15, 10).Value = 0.08182499130026955
$ws.Cells.Item(15, 12).Value = 0.3765519709796195
$ws.Cells.Item(15, 13).Value = 0.4297283454308243
$ws.Cells.Item(16, 2).Value = 1.818442938568467
$ws.Cells.Item(16, 3).Value = 0.2744432844134224
$ws.Cells.Item(16, 5).Value = 0.1695118909436371
$ws.Cells.Item(16, 6).Value = 2.517032142123
$ws.Cells.Item(16, 7).Value = 1.495417020727302
$ws.Cells.Item(16, 8).Value = 1.315168036685918
$ws.Cells.Item(16, 10).Value = 0.08222590471544322
$ws.Cells.Item(16, 12).Value = 0.3734551596131581
$ws.Cells.Item(16, 13).Value = 0.4179007565418473
$ws.Cells.Item(17, 2).Value = 1.77343509679713
$ws.Cells.Item(17, 3).Value = 0.2677206502683589
$ws.Cells.Item(17, 5).Value = 0.1696511476103977
$ws.Cells.Item(17, 6).Value = 2.512713232706332
$ws.Cells.Item(17, 7).Value = 1.489789794492424
$ws.Cells.Item(17, 8).Value = 1.31508292637514
$ws.Cells.Item(17, 10).Value = 0.08247753956053394
$ws.Cells.Item(17, 12).Value = 0.3715848081536848
$ws.Cells.Item(17, 13).Value = 0.4106844067978486
$ws.Cells.Item(18, 2).Value = 1.747605318051683
$ws.Cells.Item(18, 3).Value = 0.2638459890723084
$ws.Cells.Item(18, 5).Value = 0.1697331647610527
$ws.Cells.Item(18, 6).Value = 2.510364427610526
$ws.Cells.Item(18, 7).Value = 1.486659549034357
$ws.Cells.Item(18, 8).Value = 1.315106354875496
$ws.Cells.Item(18, 10).Value = 0.08262436568253406
$ws.Cells.Item(18, 12).Value = 0.3705198915740056
$ws.Cells.Item(18, 13).Value = 0.4065483184208105
$ws.Cells.Item(19, 2).Value = 1.738869699423333
$ws.Cells.Item(19, 3).Value = 0.2625327240773174
$ws.Cells.Item(19, 5).Value = 0.1697612645054223
$ws.Cells.Item(19, 6).Value = 2.509592389599604
$ws.Cells.Item(19, 7).Value = 1.485617942424796
$ws.Cells.Item(19, 8).Value = 1.315126708362243
$ws.Cells.Item(19, 10).Value = 0.08267443808259145
$ws.Cells.Item(19, 12).Value = 0.3701611981896207
$ws.Cells.Item(19, 13).Value = 0.4051504164773903
$ws.Cells.Item(20, 2).Value = 1.778220307296863
$ws.Cells.Item(20, 3).Value = 0.2684371118468221
$ws.Cells.Item(20, 5).Value = 0.1696361248208431
$ws.Cells.Item(20, 6).Value = 2.513158979894783
$ws.Cells.Item(20, 7).Value = 1.490377804363362
$ws.Cells.Item(20, 8).Value = 1.315084492781978
$ws.Cells.Item(20, 10).Value = 0.08245053609863717
$ws.Cells.Item(20, 12).Value = 0.3717827871897441
$ws.Cells.Item(20, 13).Value = 0.4114510932837447
$ws.Cells.Item(21, 2).Value = 1.9111191562925
$ws.Cells.Item(21, 3).Value = 0.2881797625893512
$ws.Cells.Item(21, 5).Value = 0.1692385939834979
$ws.Cells.Item(21, 6).Value = 2.526760007837098
$ws.Cells.Item(21, 7).Value = 1.507643682979278
$ws.Cells.Item(21, 8).Value = 1.315808244847261
$ws.Cells.Item(21, 10).Value = 0.08172323985583496
$ws.Cells.Item(21, 12).Value = 0.3773605764371979
$ws.Cells.Item(21, 13).Value = 0.4327941643797715
$ws.Cells.Item(22, 2).Value = 1.998484749560134
$ws.Cells.Item(22, 3).Value = 0.301011848732827
$ws.Cells.Item(22, 5).Value = 0.1689958155308267
$ws.Cells.Item(22, 6).Value = 2.536857161116458
$ws.Cells.Item(22, 7).Value = 1.519882290934078
$ws.Cells.Item(22, 8).Value = 1.316928196782072
$ws.Cells.Item(22, 10).Value = 0.08126649497850913
$ws.Cells.Item(22, 12).Value = 0.3811018405015574
$ws.Cells.Item(22, 13).Value = 0.4468717091007335
$ws.Cells.Item(23, 2).Value = 1.951809742985517
$ws.Cells.Item(23, 3).Value = 0.294169598734328
$ws.Cells.Item(23, 5).Value = 0.1691238359845384
$ws.Cells.Item(23, 6).Value = 2.531357428518334
$ws.Cells.Item(23, 7).Value = 1.513262652271834
$ws.Cells.Item(23, 8).Value = 1.316271146474634
$ws.Cells.Item(23, 10).Value = 0.08150856995280442
$ws.Cells.Item(23, 12).Value = 0.3790963275051809
$ws.Cells.Item(23, 13).Value = 0.4393465421869607
$ws.Cells.Item(24, 2).Value = 1.776056770837101
$ws.Cells.Item(24, 3).Value = 0.2681132298761781
$ws.Cells.Item(24, 5).Value = 0.1696429105291166
$ws.Cells.Item(24, 6).Value = 2.512957039626613
$ws.Cells.Item(24, 7).Value = 1.490111638358087
$ws.Cells.Item(24, 8).Value = 1.315083559239724
$ws.Cells.Item(24, 10).Value = 0.0824627376411069
$ws.Cells.Item(24, 12).Value = 0.371693248547956
$ws.Cells.Item(24, 13).Value = 0.4111044347929678
$ws.Cells.Item(25, 2).Value = 1.588680430438671
$ws.Cells.Item(25, 3).Value = 0.2397019533725882
$ws.Cells.Item(25, 5).Value = 0.1702762809285914
$ws.Cells.Item(25, 6).Value = 2.49828069481633
$ws.Cells.Item(25, 7).Value = 1.469202570734126
$ws.Cells.Item(25, 8).Value = 1.316568393789396
$ws.Cells.Item(25, 10).Value = 0.08357199495692313
$ws.Cells.Item(25, 12).Value = 0.3641235163517962
$ws.Cells.Item(25, 13).Value = 0.381198162185747
